$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) cells are stored as inline text in the source sheet (e.g. values
# like "1.609.56" that are not valid numbers, or values such as "213.57" that Excel
# would otherwise auto-convert to a number). Force each updated Price cell to Text
# format before writing the new value so it is kept as a literal string, then restore
# the cell's normal style so no stray formatting is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.789.89"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.79%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.609.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.59%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.996"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.42%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.44%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.515"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.70%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.34%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "26.90"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +11.73%  "

$ws.Range("E9").Value = "  +3.05%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0597"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.38%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0911"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.44%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.841.46"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.69%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.607.26"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.45%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.824.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.90%  "

$ws.Range("E15").Value = "  +5.54%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.57%  "

$ws.Range("E17").Value = "  +6.96%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.32"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.68%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.70%  "

$ws.Range("E20").Value = "  +3.09%  "

$ws.Range("E22").Value = "  +3.98%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.74%  "

$ws.Range("E24").Value = "  +3.61%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.03%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.108"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.44%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.39"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.56%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.37%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0472"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.22%  "

$ws.Range("E31").Value = "  +0.93%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.61%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.445.55"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.45%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.10"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.52%  "

$ws.Range("E35").Value = "  -0.21%  "

$ws.Range("E36").Value = "  +10.24%  "

$ws.Range("E37").Value = "  +2.24%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.29"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.29%  "

$ws.Range("E39").Value = "  +3.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.535"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "55.37"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +29.56%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.94"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.82%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.794"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.70%  "

$ws.Range("E44").Value = "  -0.43%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0467"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.78%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "66.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.15%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.32%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.752.07"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.89%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.87"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.20%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.835"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0103"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.82%  "
